$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("B39").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D2").Value = "64.701.84"
$ws.Range("D3").Value = "3.164.78"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "613.49"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.50"
$c.Style = "Normal"
$ws.Range("D8").Value = "3.164.56"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.528"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.47"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.474"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.79"
$c.Style = "Normal"
$ws.Range("D15").Value = "3.691.30"
$ws.Range("D17").Value = "64.698.28"
$ws.Range("D18").Value = "3.168.52"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.89"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "480.01"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.68"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.722"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.92"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "13.76"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "84.24"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.76"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.70"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "26.62"
$c.Style = "Normal"
$ws.Range("D36").Value = "0.0₃0795"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.01"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "53.33"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "462.78"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0400"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.37"
$c.Style = "Normal"
$ws.Range("D44").Value = "2.862.64"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.44"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "26.71"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "36.33"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("E29").Value = "  +3.80%  "
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("E31").Value = "  -5.56%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("E36").Value = "  +8.29%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  +3.11%  "
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +8.65%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.09%  "
